$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 122.253015
$ws.Cells.Item(2, 8).Value = 366.759045
$ws.Cells.Item(2, 9).Value = 0.1988639364328829
$ws.Cells.Item(2, 10).Value = 0.1988639364328829
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 0.8063316666666666
$ws.Cells.Item(2, 14).Value = 2.418995
$ws.Cells.Item(2, 15).Value = 0.1277387112198808
$ws.Cells.Item(2, 16).Value = 0.1277387112198808
$ws.Cells.Item(2, 17).Value = 98.57647733997499
$ws.Cells.Item(2, 18).Value = 887.188296059775
$ws.Cells.Item(2, 19).Value = 0.02540262294804876
$ws.Cells.Item(2, 20).Value = 0.02540262294804876

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 122.253015
$ws.Cells.Item(3, 8).Value = 366.759045
$ws.Cells.Item(3, 9).Value = 0.1988639364328829
$ws.Cells.Item(3, 10).Value = 0.1988639364328829
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 3.578098999999999
$ws.Cells.Item(3, 14).Value = 10.734297
$ws.Cells.Item(3, 15).Value = 0.566840884181833
$ws.Cells.Item(3, 16).Value = 0.5668408841818329
$ws.Cells.Item(3, 17).Value = 437.4333907184849
$ws.Cells.Item(3, 18).Value = 3936.900516466364
$ws.Cells.Item(3, 19).Value = 0.1127242095594952
$ws.Cells.Item(3, 20).Value = 0.1127242095594951

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 122.253015
$ws.Cells.Item(4, 8).Value = 366.759045
$ws.Cells.Item(4, 9).Value = 0.1988639364328829
$ws.Cells.Item(4, 10).Value = 0.1988639364328829
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.215895
$ws.Cells.Item(4, 14).Value = 0.647685
$ws.Cells.Item(4, 15).Value = 0.03420199180918047
$ws.Cells.Item(4, 16).Value = 0.03420199180918047
$ws.Cells.Item(4, 17).Value = 26.393814673425
$ws.Cells.Item(4, 18).Value = 237.544332060825
$ws.Cells.Item(4, 19).Value = 0.006801542725018845
$ws.Cells.Item(4, 20).Value = 0.006801542725018844

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 122.253015
$ws.Cells.Item(5, 8).Value = 366.759045
$ws.Cells.Item(5, 9).Value = 0.1988639364328829
$ws.Cells.Item(5, 10).Value = 0.1988639364328829
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.712026
$ws.Cells.Item(5, 14).Value = 5.136078
$ws.Cells.Item(5, 15).Value = 0.2712184127891059
$ws.Cells.Item(5, 16).Value = 0.2712184127891059
$ws.Cells.Item(5, 17).Value = 209.30034025839
$ws.Cells.Item(5, 18).Value = 1883.70306232551
$ws.Cells.Item(5, 19).Value = 0.05393556120032014
$ws.Cells.Item(5, 20).Value = 0.05393556120032013

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 132.5447616666667
$ws.Cells.Item(6, 8).Value = 397.634285
$ws.Cells.Item(6, 9).Value = 0.2156050961899926
$ws.Cells.Item(6, 10).Value = 0.2156050961899926
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.8063316666666666
$ws.Cells.Item(6, 14).Value = 2.418995
$ws.Cells.Item(6, 15).Value = 0.1277387112198808
$ws.Cells.Item(6, 16).Value = 0.1277387112198808
$ws.Cells.Item(6, 17).Value = 106.8750385826194
$ws.Cells.Item(6, 18).Value = 961.8753472435749
$ws.Cells.Item(6, 19).Value = 0.02754111711974809
$ws.Cells.Item(6, 20).Value = 0.02754111711974809

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 132.5447616666667
$ws.Cells.Item(7, 8).Value = 397.634285
$ws.Cells.Item(7, 9).Value = 0.2156050961899926
$ws.Cells.Item(7, 10).Value = 0.2156050961899926
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 3.578098999999999
$ws.Cells.Item(7, 14).Value = 10.734297
$ws.Cells.Item(7, 15).Value = 0.566840884181833
$ws.Cells.Item(7, 16).Value = 0.5668408841818329
$ws.Cells.Item(7, 17).Value = 474.2582791747382
$ws.Cells.Item(7, 18).Value = 4268.324512572644
$ws.Cells.Item(7, 19).Value = 0.1222137833584446
$ws.Cells.Item(7, 20).Value = 0.1222137833584445

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 132.5447616666667
$ws.Cells.Item(8, 8).Value = 397.634285
$ws.Cells.Item(8, 9).Value = 0.2156050961899926
$ws.Cells.Item(8, 10).Value = 0.2156050961899926
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 0.215895
$ws.Cells.Item(8, 14).Value = 0.647685
$ws.Cells.Item(8, 15).Value = 0.03420199180918047
$ws.Cells.Item(8, 16).Value = 0.03420199180918047
$ws.Cells.Item(8, 17).Value = 28.61575132002499
$ws.Cells.Item(8, 18).Value = 257.541761880225
$ws.Cells.Item(8, 19).Value = 0.007374123733907694
$ws.Cells.Item(8, 20).Value = 0.007374123733907694

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 132.5447616666667
$ws.Cells.Item(9, 8).Value = 397.634285
$ws.Cells.Item(9, 9).Value = 0.2156050961899926
$ws.Cells.Item(9, 10).Value = 0.2156050961899926
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 1.712026
$ws.Cells.Item(9, 14).Value = 5.136078
$ws.Cells.Item(9, 15).Value = 0.2712184127891059
$ws.Cells.Item(9, 16).Value = 0.2712184127891059
$ws.Cells.Item(9, 17).Value = 226.9200781371367
$ws.Cells.Item(9, 18).Value = 2042.28070323423
$ws.Cells.Item(9, 19).Value = 0.05847607197789229
$ws.Cells.Item(9, 20).Value = 0.05847607197789229

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 320.0894206666666
$ws.Cells.Item(10, 8).Value = 960.2682619999999
$ws.Cells.Item(10, 9).Value = 0.5206762565675317
$ws.Cells.Item(10, 10).Value = 0.5206762565675317
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 0.8063316666666666
$ws.Cells.Item(10, 14).Value = 2.418995
$ws.Cells.Item(10, 15).Value = 0.1277387112198808
$ws.Cells.Item(10, 16).Value = 0.1277387112198808
$ws.Cells.Item(10, 17).Value = 258.098236048521
$ws.Cells.Item(10, 18).Value = 2322.884124436689
$ws.Cells.Item(10, 19).Value = 0.0665105139767285
$ws.Cells.Item(10, 20).Value = 0.0665105139767285

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 320.0894206666666
$ws.Cells.Item(11, 8).Value = 960.2682619999999
$ws.Cells.Item(11, 9).Value = 0.5206762565675317
$ws.Cells.Item(11, 10).Value = 0.5206762565675317
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 3.578098999999999
$ws.Cells.Item(11, 14).Value = 10.734297
$ws.Cells.Item(11, 15).Value = 0.566840884181833
$ws.Cells.Item(11, 16).Value = 0.5668408841818329
$ws.Cells.Item(11, 17).Value = 1145.311635997979
$ws.Cells.Item(11, 18).Value = 10307.80472398181
$ws.Cells.Item(11, 19).Value = 0.2951405896452266
$ws.Cells.Item(11, 20).Value = 0.2951405896452265

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 320.0894206666666
$ws.Cells.Item(12, 8).Value = 960.2682619999999
$ws.Cells.Item(12, 9).Value = 0.5206762565675317
$ws.Cells.Item(12, 10).Value = 0.5206762565675317
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.215895
$ws.Cells.Item(12, 14).Value = 0.647685
$ws.Cells.Item(12, 15).Value = 0.03420199180918047
$ws.Cells.Item(12, 16).Value = 0.03420199180918047
$ws.Cells.Item(12, 17).Value = 69.10570547482999
$ws.Cells.Item(12, 18).Value = 621.9513492734699
$ws.Cells.Item(12, 19).Value = 0.01780816506235747
$ws.Cells.Item(12, 20).Value = 0.01780816506235747

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 320.0894206666666
$ws.Cells.Item(13, 8).Value = 960.2682619999999
$ws.Cells.Item(13, 9).Value = 0.5206762565675317
$ws.Cells.Item(13, 10).Value = 0.5206762565675317
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 1.712026
$ws.Cells.Item(13, 14).Value = 5.136078
$ws.Cells.Item(13, 15).Value = 0.2712184127891059
$ws.Cells.Item(13, 16).Value = 0.2712184127891059
$ws.Cells.Item(13, 17).Value = 548.0014105062706
$ws.Cells.Item(13, 18).Value = 4932.012694556436
$ws.Cells.Item(13, 19).Value = 0.1412169878832192
$ws.Cells.Item(13, 20).Value = 0.1412169878832192

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 39.86989333333333
$ws.Cells.Item(14, 8).Value = 119.60968
$ws.Cells.Item(14, 9).Value = 0.06485471080959287
$ws.Cells.Item(14, 10).Value = 0.06485471080959287
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 0.8063316666666666
$ws.Cells.Item(14, 14).Value = 2.418995
$ws.Cells.Item(14, 15).Value = 0.1277387112198808
$ws.Cells.Item(14, 16).Value = 0.1277387112198808
$ws.Cells.Item(14, 17).Value = 32.14835754128888
$ws.Cells.Item(14, 18).Value = 289.3352178716
$ws.Cells.Item(14, 19).Value = 0.008284457175355467
$ws.Cells.Item(14, 20).Value = 0.008284457175355467

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 39.86989333333333
$ws.Cells.Item(15, 8).Value = 119.60968
$ws.Cells.Item(15, 9).Value = 0.06485471080959287
$ws.Cells.Item(15, 10).Value = 0.06485471080959287
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 3.578098999999999
$ws.Cells.Item(15, 14).Value = 10.734297
$ws.Cells.Item(15, 15).Value = 0.566840884181833
$ws.Cells.Item(15, 16).Value = 0.5668408841818329
$ws.Cells.Item(15, 17).Value = 142.6584254661066
$ws.Cells.Item(15, 18).Value = 1283.92582919496
$ws.Cells.Item(15, 19).Value = 0.03676230161866671
$ws.Cells.Item(15, 20).Value = 0.0367623016186667

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 39.86989333333333
$ws.Cells.Item(16, 8).Value = 119.60968
$ws.Cells.Item(16, 9).Value = 0.06485471080959287
$ws.Cells.Item(16, 10).Value = 0.06485471080959287
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.215895
$ws.Cells.Item(16, 14).Value = 0.647685
$ws.Cells.Item(16, 15).Value = 0.03420199180918047
$ws.Cells.Item(16, 16).Value = 0.03420199180918047
$ws.Cells.Item(16, 17).Value = 8.607710621199999
$ws.Cells.Item(16, 18).Value = 77.46939559079999
$ws.Cells.Item(16, 19).Value = 0.002218160287896463
$ws.Cells.Item(16, 20).Value = 0.002218160287896463

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 39.86989333333333
$ws.Cells.Item(17, 8).Value = 119.60968
$ws.Cells.Item(17, 9).Value = 0.06485471080959287
$ws.Cells.Item(17, 10).Value = 0.06485471080959287
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 1.712026
$ws.Cells.Item(17, 14).Value = 5.136078
$ws.Cells.Item(17, 15).Value = 0.2712184127891059
$ws.Cells.Item(17, 16).Value = 0.2712184127891059
$ws.Cells.Item(17, 17).Value = 68.25829400389333
$ws.Cells.Item(17, 18).Value = 614.32464603504
$ws.Cells.Item(17, 19).Value = 0.01758979172767425
$ws.Cells.Item(17, 20).Value = 0.01758979172767425
